$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '38.378.93'
$ws.Range("D2").NumberFormat = "General"
$ws.Range("E2").Value = '  +1.55%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.088.27'
$ws.Range("D3").NumberFormat = "General"
$ws.Range("E3").Value = '  +2.29%  '

$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("E5").Value = '  -0.22%  '

$ws.Range("E6").Value = '  +0.70%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '60.75'
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = '  -0.30%  '

$ws.Range("E9").Value = '  +1.43%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0835'
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = '  +1.56%  '

$ws.Range("E11").Value = '  -0.20%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '2.398.47'
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = '  +2.28%  '

$ws.Range("E13").Value = '  +1.46%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '22.32'
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = '  +6.24%  '

$ws.Range("E15").Value = '  +0.97%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.43'
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = '  +4.46%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.093.76'
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = '  +2.28%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '38.285.42'
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = '  +1.26%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '71.09'
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = '  +2.22%  '

$ws.Range("E20").Value = '  +2.78%  '

$ws.Range("E21").Value = '  +0.96%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '225.32'
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = '  +0.69%  '

$ws.Range("E23").Value = '  +0.06%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.44'
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = '  +0.45%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.31'
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = '  +1.87%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '169.61'
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = '  +0.92%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.42'
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = '  +0.93%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.134'
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = '  +3.85%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.04'
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = '  +1.27%  '

$ws.Range("E30").Value = '  +8.18%  '

$ws.Range("E31").Value = '  -0.41%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.35'
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = '  +5.92%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.79'
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = '  +6.31%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.50'
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = '  +2.80%  '

$ws.Range("E35").Value = '  +0.59%  '

$ws.Range("B36").Value = 'THORChain'
$ws.Range("C36").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.43'
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = '  -2.49%  '

$ws.Range("B37").Value = 'LidoDAOToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.38'
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = '  +2.19%  '

$ws.Range("E38").Value = '  +3.09%  '

$ws.Range("E39").Value = '  -0.17%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '18.51'
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = '  +2.44%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.538.50'
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = '  +0.03%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '99.97'
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = '  +4.02%  '

$ws.Range("E43").Value = '  +1.31%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0932'
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = '  +2.30%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.81'
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = '  -0.18%  '

$ws.Range("E46").Value = '  +10.08%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.14'
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = '  -1.47%  '

$ws.Range("E48").Value = '  +0.89%  '

$ws.Range("E49").Value = '  +2.48%  '

$ws.Range("E50").Value = '  +1.10%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.283.72'
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = '  +2.19%  '
